$d = $word.ActiveDocument

function New-PkgXml($innerBodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Step 1: strip the _GoBack bookmark off the end of the "...Fancybox for
# lightboxes." (Technology) paragraph, leaving its runs completely
# untouched (just delete the bookmark itself, via the named-bookmark
# lookup, rather than rewriting the paragraph). ---
$fancyPara = $d.Paragraphs(38)
if ($fancyPara.Range.Text -notlike "*Fancybox for lightboxes.*") {
    throw "Paragraph 38 is not the expected Technology/Fancybox paragraph: $($fancyPara.Range.Text)"
}
$goBack = $d.Bookmarks("_GoBack")
if ($goBack.Range.Text -notlike "*Fancybox for lightboxes.*") {
    throw "_GoBack bookmark is not where expected: $($goBack.Range.Text)"
}
$goBack.Delete()

# --- Step 2: replace the "Knowledge" section body paragraph (the old
# roster/API text) with the new lightbox text, and move the _GoBack
# bookmark here (at the very end of the paragraph). ---
$knowledgePara = $d.Paragraphs(44)
if ($knowledgePara.Range.Text -notlike "*roster page working*") {
    throw "Paragraph 44 is not the expected Knowledge-body paragraph: $($knowledgePara.Range.Text)"
}
$knowledgeXml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I need to find an easy-to-use lightbox to use for the News page (in order to view the images i</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>n a larger scale).</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$knowledgePara.Range.InsertXML((New-PkgXml $knowledgeXml))

# --- Step 3: delete the trailing empty-ish paragraph that only held a
# single space character. ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
if ($lastPara.Range.Text.TrimEnd("`r") -ne " ") {
    throw "Last paragraph is not the expected single-space paragraph: [$($lastPara.Range.Text)]"
}
$lastPara.Range.Delete()

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
